$d = $word.ActiveDocument

function Replace-ParagraphXml($para, $bodyXml) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($xml)
}

# --- Change 1: split "If somebody else cloned..." paragraph (which holds the _GoBack
# bookmark) into five paragraphs, adding two new paragraphs of text and moving the
# bookmark into its own empty paragraph. ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "If somebody else cloned my repository*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    Write-Host "ERROR: could not find paragraph 'If somebody else cloned my repository...'"
}

$body1 = @'
<w:p><w:pPr><w:pStyle w:val="Standard"/></w:pPr><w:r><w:t>If somebody else cloned my repository to their local computer, they could clone the new file using ‘git pull’. This pulls down from GitHub to their computer, however, it only works if they’re in the correct repo/ working directory.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Standard"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Standard"/></w:pPr><w:r><w:t>If we make changes to the file and then look at ‘git status’, we can see that the file has been modified. To update the new file, and push it to the glo</w:t></w:r><w:r><w:t>bal repo: first type ‘git commit –m “instructions on pushing modified files”</w:t></w:r><w:r><w:t xml:space="preserve">’, then ‘git push’. Check the file is now on GitHub.com. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Standard"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="Standard"/></w:pPr><w:r><w:t>To clone a pre-made repository from GutHub: i.e. clone the repository for the Statistical Inference John Hopkins coursera course. First create a repo on GitHub, termed Statistical-Inference-JH and initialise with a readme file. Then clone the repo to my computer as before, so we have a new directory called ‘Statistical-Inference-JH’ under within the ‘GitHub’ directory. Then type ‘git clone  https://github.com/bcaffo/courses.git’ (which is the URL for the entire JH repo). These files will now be saved locally to my computer.</w:t></w:r></w:p>
'@

if ($target -ne $null) {
    Replace-ParagraphXml $target $body1
}

# --- Change 2: remove the stray <w:lastRenderedPageBreak/> before the "ls " run. ---
$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "ls *: lists the files in the wd*") {
        $target2 = $p
        break
    }
}
if ($target2 -eq $null) {
    Write-Host "ERROR: could not find paragraph 'ls : lists the files in the wd'"
}

$body2 = @'
<w:p><w:pPr><w:pStyle w:val="Standard"/></w:pPr><w:r><w:t xml:space="preserve">ls </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>: lists the files in the wd</w:t></w:r></w:p>
'@

if ($target2 -ne $null) {
    Replace-ParagraphXml $target2 $body2
}

Write-Host "Edits applied"
